$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.918.01'
$ws.Range("E2").Value = '  +0.04%  '

$ws.Range("D3").Value = '3.555.09'
$ws.Range("E3").Value = '  +2.49%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("D7").Value = '3.555.24'
$ws.Range("E7").Value = '  +2.45%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +0.52%  '

$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("E11").Value = '  -2.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.384'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.81%  '

$ws.Range("D13").Value = '4.160.78'
$ws.Range("E13").Value = '  +2.55%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000181'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.32%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.557.90'
$ws.Range("E15").Value = '  +2.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.87'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '65.022.30'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("E23").Value = '  +3.68%  '

$ws.Range("D24").Value = '3.701.35'
$ws.Range("E24").Value = '  +2.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.38%  '

$ws.Range("E27").Value = '  +4.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.37%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +3.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.22%  '

$ws.Range("E32").Value = '  +25.61%  '

$ws.Range("D33").Value = '3.554.57'
$ws.Range("E33").Value = '  +1.96%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.44%  '

$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.144'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.30%  '

$ws.Range("E40").Value = '  +5.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0803'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.66%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.08%  '

$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("E46").Value = '  +2.95%  '

$ws.Range("E47").Value = '  +4.48%  '

$ws.Range("E48").Value = '  +2.07%  '

$ws.Range("D49").Value = '2.482.81'
$ws.Range("E49").Value = '  +12.20%  '

$ws.Range("E50").Value = '  +3.58%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0260'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.67%  '
